$wb = $excel.ActiveWorkbook

# ===== Sheet 1: LP1912 =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 18:44:57"
$ws.Range("A3").Value = "Total filas: 322"
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "07:50:33"
$arr[0,1] = "09:39"
$arr[0,2] = "15_ABASTO"
$arr[0,3] = 109
$arr[0,4] = "LP1912"
$ws.Range("A64:E64").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "08:27:16"
$arr[0,1] = "09:39"
$arr[0,2] = "23_HERNANDEZ"
$arr[0,3] = 72
$arr[0,4] = "LP1912"
$ws.Range("A65:E65").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:50:41"
$arr[0,1] = "12:36"
$arr[0,2] = "27_EL RETIRO"
$arr[0,3] = 106
$arr[0,4] = "LP1912"
$ws.Range("A135:E135").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:34:59"
$arr[0,1] = "12:36"
$arr[0,2] = "23_HERNANDEZ"
$arr[0,3] = 62
$arr[0,4] = "LP1912"
$ws.Range("A136:E136").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "12:33:21"
$arr[0,1] = "14:17"
$arr[0,2] = "27_EL RETIRO"
$arr[0,3] = 104
$arr[0,4] = "LP1912"
$ws.Range("A176:E176").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "12:33:21"
$arr[0,1] = "14:17"
$arr[0,2] = "11_ETCHEVERRY"
$arr[0,3] = 104
$arr[0,4] = "LP1912"
$ws.Range("A177:E177").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "15:17:33"
$arr[0,1] = "16:05"
$arr[0,2] = "16_SANTA ANA"
$arr[0,3] = 48
$arr[0,4] = "LP1912"
$ws.Range("A216:E216").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "14:12:26"
$arr[0,1] = "16:05"
$arr[0,2] = "14_ABASTO"
$arr[0,3] = 113
$arr[0,4] = "LP1912"
$ws.Range("A217:E217").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "15:58:05"
$arr[0,1] = "17:46"
$arr[0,2] = "215_EL PELIGRO"
$arr[0,3] = 108
$arr[0,4] = "LP1912"
$ws.Range("A257:E257").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "16:52:42"
$arr[0,1] = "17:46"
$arr[0,2] = "23_HERNANDEZ"
$arr[0,3] = 54
$arr[0,4] = "LP1912"
$ws.Range("A258:E258").Value = $arr
$ws.Rows("312:327").Insert()
$arr = New-Object 'object[,]' 39,5
$arr[0,0] = "18:44:57"
$arr[0,1] = "18:51"
$arr[0,2] = "17_ROMERO"
$arr[0,3] = 7
$arr[0,4] = "LP1912"
$arr[1,0] = "17:13:39"
$arr[1,1] = "18:52"
$arr[1,2] = "17_ROMERO"
$arr[1,3] = 99
$arr[1,4] = "LP1912"
$arr[2,0] = "18:44:57"
$arr[2,1] = "18:56"
$arr[2,2] = "16_P MOR-SANTA ANA"
$arr[2,3] = 12
$arr[2,4] = "LP1912"
$arr[3,0] = "17:13:39"
$arr[3,1] = "18:57"
$arr[3,2] = "16_P MOR-SANTA ANA"
$arr[3,3] = 104
$arr[3,4] = "LP1912"
$arr[4,0] = "17:13:39"
$arr[4,1] = "18:59"
$arr[4,2] = "14_ABASTO"
$arr[4,3] = 106
$arr[4,4] = "LP1912"
$arr[5,0] = "17:36:10"
$arr[5,1] = "19:00"
$arr[5,2] = "14_ABASTO"
$arr[5,3] = 84
$arr[5,4] = "LP1912"
$arr[6,0] = "17:13:39"
$arr[6,1] = "19:03"
$arr[6,2] = "215_EL PELIGRO"
$arr[6,3] = 110
$arr[6,4] = "LP1912"
$arr[7,0] = "17:36:10"
$arr[7,1] = "19:04"
$arr[7,2] = "215_EL PELIGRO"
$arr[7,3] = 88
$arr[7,4] = "LP1912"
$arr[8,0] = "17:56:03"
$arr[8,1] = "19:10"
$arr[8,2] = "27_EL RETIRO"
$arr[8,3] = 74
$arr[8,4] = "LP1912"
$arr[9,0] = "18:12:30"
$arr[9,1] = "19:10"
$arr[9,2] = "16_SANTA ANA"
$arr[9,3] = 58
$arr[9,4] = "LP1912"
$arr[10,0] = "17:48:33"
$arr[10,1] = "19:12"
$arr[10,2] = "27_EL RETIRO"
$arr[10,3] = 84
$arr[10,4] = "LP1912"
$arr[11,0] = "18:44:57"
$arr[11,1] = "19:15"
$arr[11,2] = "17_ROMERO"
$arr[11,3] = 31
$arr[11,4] = "LP1912"
$arr[12,0] = "18:12:30"
$arr[12,1] = "19:16"
$arr[12,2] = "27_EL RETIRO"
$arr[12,3] = 64
$arr[12,4] = "LP1912"
$arr[13,0] = "18:44:57"
$arr[13,1] = "19:16"
$arr[13,2] = "14X44_ABASTO"
$arr[13,3] = 32
$arr[13,4] = "LP1912"
$arr[14,0] = "17:56:03"
$arr[14,1] = "19:16"
$arr[14,2] = "17_ROMERO"
$arr[14,3] = 80
$arr[14,4] = "LP1912"
$arr[15,0] = "17:36:10"
$arr[15,1] = "19:17"
$arr[15,2] = "27_EL RETIRO"
$arr[15,3] = 101
$arr[15,4] = "LP1912"
$arr[16,0] = "17:36:10"
$arr[16,1] = "19:17"
$arr[16,2] = "14X44_ABASTO"
$arr[16,3] = 101
$arr[16,4] = "LP1912"
$arr[17,0] = "17:56:03"
$arr[17,1] = "19:21"
$arr[17,2] = "23_HERNANDEZ"
$arr[17,3] = 85
$arr[17,4] = "LP1912"
$arr[18,0] = "18:12:30"
$arr[18,1] = "19:22"
$arr[18,2] = "23_HERNANDEZ"
$arr[18,3] = 70
$arr[18,4] = "LP1912"
$arr[19,0] = "18:44:57"
$arr[19,1] = "19:23"
$arr[19,2] = "16_SANTA ANA"
$arr[19,3] = 39
$arr[19,4] = "LP1912"
$arr[20,0] = "18:44:57"
$arr[20,1] = "19:27"
$arr[20,2] = "215C_EL PATO"
$arr[20,3] = 43
$arr[20,4] = "LP1912"
$arr[21,0] = "18:44:57"
$arr[21,1] = "19:27"
$arr[21,2] = "16_P MOR-SANTA ANA"
$arr[21,3] = 43
$arr[21,4] = "LP1912"
$arr[22,0] = "17:36:10"
$arr[22,1] = "19:28"
$arr[22,2] = "215C_EL PATO"
$arr[22,3] = 112
$arr[22,4] = "LP1912"
$arr[23,0] = "17:48:33"
$arr[23,1] = "19:35"
$arr[23,2] = "11_ETCHEVERRY"
$arr[23,3] = 107
$arr[23,4] = "LP1912"
$arr[24,0] = "17:56:03"
$arr[24,1] = "19:36"
$arr[24,2] = "11_ETCHEVERRY"
$arr[24,3] = 100
$arr[24,4] = "LP1912"
$arr[25,0] = "18:44:57"
$arr[25,1] = "19:38"
$arr[25,2] = "15X38_ABASTO"
$arr[25,3] = 54
$arr[25,4] = "LP1912"
$arr[26,0] = "17:48:33"
$arr[26,1] = "19:39"
$arr[26,2] = "15X38_ABASTO"
$arr[26,3] = 111
$arr[26,4] = "LP1912"
$arr[27,0] = "18:44:57"
$arr[27,1] = "19:51"
$arr[27,2] = "81_EL PELIGRO"
$arr[27,3] = 67
$arr[27,4] = "LP1912"
$arr[28,0] = "17:56:03"
$arr[28,1] = "19:52"
$arr[28,2] = "81_EL PELIGRO"
$arr[28,3] = 116
$arr[28,4] = "LP1912"
$arr[29,0] = "18:44:57"
$arr[29,1] = "19:52"
$arr[29,2] = "225_GOMEZ"
$arr[29,3] = 68
$arr[29,4] = "LP1912"
$arr[30,0] = "17:56:03"
$arr[30,1] = "19:53"
$arr[30,2] = "225_GOMEZ"
$arr[30,3] = 117
$arr[30,4] = "LP1912"
$arr[31,0] = "18:44:57"
$arr[31,1] = "19:53"
$arr[31,2] = "16_SANTA ANA"
$arr[31,3] = 69
$arr[31,4] = "LP1912"
$arr[32,0] = "18:44:57"
$arr[32,1] = "20:06"
$arr[32,2] = "215C_EL PATO"
$arr[32,3] = 82
$arr[32,4] = "LP1912"
$arr[33,0] = "18:44:57"
$arr[33,1] = "20:09"
$arr[33,2] = "23_HERNANDEZ"
$arr[33,3] = 85
$arr[33,4] = "LP1912"
$arr[34,0] = "18:31:25"
$arr[34,1] = "20:12"
$arr[34,2] = "215C_EL PATO"
$arr[34,3] = 101
$arr[34,4] = "LP1912"
$arr[35,0] = "18:44:57"
$arr[35,1] = "20:12"
$arr[35,2] = "14_ABASTO"
$arr[35,3] = 88
$arr[35,4] = "LP1912"
$arr[36,0] = "18:44:57"
$arr[36,1] = "20:21"
$arr[36,2] = "15_ABASTO"
$arr[36,3] = 97
$arr[36,4] = "LP1912"
$arr[37,0] = "18:31:25"
$arr[37,1] = "20:22"
$arr[37,2] = "15_ABASTO"
$arr[37,3] = 111
$arr[37,4] = "LP1912"
$arr[38,0] = "18:44:57"
$arr[38,1] = "20:30"
$arr[38,2] = "10_OLMOS"
$arr[38,3] = 106
$arr[38,4] = "LP1912"
$ws.Range("A289:E327").Value = $arr

# ===== Sheet 2: LP1912-215 =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 18:44:57"
$ws.Range("A3").Value = "Total filas: 51"
$ws.Rows("55:56").Insert()
$arr = New-Object 'object[,]' 4,5
$arr[0,0] = "18:44:57"
$arr[0,1] = "19:27"
$arr[0,2] = "215C_EL PATO"
$arr[0,3] = 43
$arr[0,4] = "LP1912"
$arr[1,0] = "17:36:10"
$arr[1,1] = "19:28"
$arr[1,2] = "215C_EL PATO"
$arr[1,3] = 112
$arr[1,4] = "LP1912"
$arr[2,0] = "18:44:57"
$arr[2,1] = "20:06"
$arr[2,2] = "215C_EL PATO"
$arr[2,3] = 82
$arr[2,4] = "LP1912"
$arr[3,0] = "18:31:25"
$arr[3,1] = "20:12"
$arr[3,2] = "215C_EL PATO"
$arr[3,3] = 101
$arr[3,4] = "LP1912"
$ws.Range("A53:E56").Value = $arr

# ===== Sheet 3: 6203-6173 =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 18:44:57"
$ws.Range("A3").Value = "Total filas: 41"
$ws.Rows("45:46").Insert()
$arr = New-Object 'object[,]' 4,5
$arr[0,0] = "18:44:57"
$arr[0,1] = "19:23"
$arr[0,2] = "215B_LP-P MOR-1 Y 57"
$arr[0,3] = 39
$arr[0,4] = "L6173"
$arr[1,0] = "17:36:10"
$arr[1,1] = "19:24"
$arr[1,2] = "215B_LP-P MOR-1 Y 57"
$arr[1,3] = 108
$arr[1,4] = "L6173"
$arr[2,0] = "18:44:57"
$arr[2,1] = "19:57"
$arr[2,2] = "215C_LA PLATA"
$arr[2,3] = 73
$arr[2,4] = "L6203"
$arr[3,0] = "18:12:30"
$arr[3,1] = "19:58"
$arr[3,2] = "215C_LA PLATA"
$arr[3,3] = 106
$arr[3,4] = "L6203"
$ws.Range("A43:E46").Value = $arr

